# This workbook stores a weekly price table for "Tuna" (prickly pear) at the
# "Agricola del Norte S.A. de Arica" market. The edit re-assigns the weekly
# observation data (Fecha, Calidad, Volumen, Precios, Unidad, Origen, etc.) to
# different rows -- i.e. it is a permutation of the data across rows 2-18
# (row 5 keeps its original data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D","L","M","N","O","P","Q","R","S","T")

# Snapshot the current values of the affected columns for every data row
# before overwriting anything, since several rows trade values with each other.
$snapshot = @{}
foreach ($r in 2..18) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Maps each target row to the row whose original data it should now hold.
$mapping = @{
    2 = 3
    3 = 13
    4 = 14
    6 = 9
    7 = 10
    8 = 16
    9 = 15
    10 = 17
    11 = 18
    12 = 7
    13 = 8
    14 = 12
    15 = 6
    16 = 2
    17 = 4
    18 = 11
}

foreach ($target in $mapping.Keys) {
    $source = $mapping[$target]
    $srcVals = $snapshot[$source]
    foreach ($c in $cols) {
        $ws.Range("$c$target").Value = $srcVals[$c]
    }
}
